# Apply the commit: "Fixed model_id in handle_info_request function and
# updated differences_data filtering"
#
# 1) Sheet "اطلاعات کلی" (general info): column F ("ویژگی ID" / feature id)
#    is re-mapped so several device variants now share the same feature id
#    (the "model_id" bug fix).
# 2) Sheet "ویژگی" (features): the feature text for each row gets a
#    device-specific suffix appended so every feature id maps to a unique,
#    more descriptive string.
# 3) Sheet "معرفی" (intro): the intro text for each row gets a
#    device-specific suffix appended, same idea as above.
# 4) The active/selected tab moves from "اطلاعات کلی" (sheet 1) to
#    "تفاوت" (sheet 5).

$wb = $excel.ActiveWorkbook

# --- 1) اطلاعات کلی : fix column F (ویژگی ID) ------------------------------
$wsInfo = $wb.Worksheets.Item("اطلاعات کلی")

$wsInfo.Range("F3").Value = 1
$wsInfo.Range("F4").Value = 1

$wsInfo.Range("F5").Value = 2
$wsInfo.Range("F6").Value = 2
$wsInfo.Range("F7").Value = 2

$wsInfo.Range("F8").Value = 3
$wsInfo.Range("F9").Value = 3
$wsInfo.Range("F10").Value = 3

$wsInfo.Range("F11").Value = 4
$wsInfo.Range("F12").Value = 4
$wsInfo.Range("F13").Value = 4

$wsInfo.Range("F14").Value = 5
$wsInfo.Range("F15").Value = 5

# --- 2) ویژگی : append device-specific detail to each feature row --------
$wsFeature = $wb.Worksheets.Item("ویژگی")

$wsFeature.Range("B2").Value  = "دقت بالا , حکاکی سریع (موپا دستی هوشمند)"
$wsFeature.Range("B3").Value  = "دقت بالا , حکاکی سریع (موپا میزدار هوشمند)"
$wsFeature.Range("B4").Value  = "دقت بالا , حکاکی سریع (موپا روی میز هوشمند)"
$wsFeature.Range("B5").Value  = "قابلیت تنظیم هوشمند,سرعت پردازش بالا (فایبر دستی هوشمند)"
$wsFeature.Range("B6").Value  = "قابلیت تنظیم هوشمند,سرعت پردازش بالا (فایبر میزدار هوشمند)"
$wsFeature.Range("B7").Value  = "قابلیت تنظیم هوشمند,سرعت پردازش بالا (فایبر وی میز هوشمند)"
$wsFeature.Range("B8").Value  = "مناسب برای مواد غیر فلزی (فایبر دستی غیر هوشمند)"
$wsFeature.Range("B9").Value  = "مناسب برای مواد غیر فلزی (فایبر میزدار غیر هوشمند)"
$wsFeature.Range("B10").Value = "مناسب برای مواد غیر فلزی (فایبر روی میز غیر هوشمند)"
$wsFeature.Range("B11").Value = "سرعت برش بسیار بالا (40wneje)"
$wsFeature.Range("B12").Value = "سرعت برش بسیار بالا(neje80w)"
$wsFeature.Range("B13").Value = "سرعت برش بسیار بالا (neje160w)"
$wsFeature.Range("B14").Value = "باکس صنعتی و مناسب برای محیط های صنعتی(neje200 w)"
$wsFeature.Range("B15").Value = "باکس صنعتی و مناسب برای محیط های صنعتی(neje200w باکس دار)"

# --- 3) معرفی : append device-specific detail to each intro row ----------
$wsIntro = $wb.Worksheets.Item("معرفی")

$wsIntro.Range("B2").Value  = "دستگاه موپا برا یحکاکی دقیق روی فلزات (موپا  دستی هوشمند)"
$wsIntro.Range("B3").Value  = "دستگاه موپا برا یحکاکی دقیق روی فلزات (موپا میزدار هوشمند)"
$wsIntro.Range("B4").Value  = "دستگاه موپا برا یحکاکی دقیق روی فلزات (موپا روی میز هوشمند)"
$wsIntro.Range("B5").Value  = "فایبر مناسب برای فلزات سنگین و متنوع (فایبر دستی هوشمند)"
$wsIntro.Range("B6").Value  = "فایبر مناسب برای فلزات سنگین و متنوع (فایبر میزدار هوشمند)"
$wsIntro.Range("B7").Value  = "فایبر مناسب برای فلزات سنگین و متنوع (فایبر روی میز هوشمند)"
$wsIntro.Range("B8").Value  = "فایبر مناسب برای فلزات سنگین و متنوع (فایبر دستی غیر هشمند)"
$wsIntro.Range("B9").Value  = "فایبر مناسب برای فلزات سنگین و متنوع (فایبر میزدار غیر هوشمند)"
$wsIntro.Range("B10").Value = "فایبر مناسب برای فلزات سنگین و متنوع ( فایبر روی میز غیر هوشمند)"
$wsIntro.Range("B11").Value = "دیود برای مواد غیر فلزی مثل چوب و شیشه مناسب است (neje40w)"
$wsIntro.Range("B12").Value = "دیود برای مواد غیر فلزی مثل چوب و شیشه مناسب است (neje80w)"
$wsIntro.Range("B13").Value = "دیود برای مواد غیر فلزی مثل چوب و شیشه مناسب است (neje160w)"
$wsIntro.Range("B14").Value = "دیود برای مواد غیر فلزی مثل چوب و شیشه مناسب است(neje200w هوشمند)"
$wsIntro.Range("B15").Value = "دیود برای مواد غیر فلزی مثل چوب و شیشه مناسب است(neje200w باکس دار)"

# --- 4) Move the active tab from "اطلاعات کلی" to "تفاوت" -----------------
$wsDiff = $wb.Worksheets.Item("تفاوت")
$wsDiff.Activate()
